# daily auto push: 2025-09-30 22:29 UTC
# Append the new daily-ranking record as row 43 (date 2025/10/01, weekday 水,
# time 6, ranking 154), extending the sheet's used range from A1:D42 to A1:D43.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds dates stored as plain text (e.g. "2025/09/30"), not real
# Excel date serials. Force the cell to Text format before assigning the
# value so the "YYYY/MM/DD" string isn't auto-converted into a date, then
# clear the formatting again so the new cell ends up unstyled, matching the
# rest of the data rows.
$ws.Range("A43").NumberFormat = "@"
$ws.Range("A43").Value = "2025/10/01"
$ws.Range("A43").ClearFormats()

$ws.Range("B43").Value = "水"
$ws.Range("C43").Value = 6
$ws.Range("D43").Value = 154
